$wb = $excel.ActiveWorkbook

# Users sheet: Gabriella Ghosn -> Amanda Donovan
$usersWs = $wb.Worksheets.Item("Users")
$usersWs.Range("A2").Value = "Amanda Donovan"

# MoreAttendees sheet: James Craven -> Thomas Bailey
$moreWs = $wb.Worksheets.Item("MoreAttendees")
$moreWs.Range("B2").Value = "Thomas Bailey"
$moreWs.Range("J21").Select() | Out-Null

# Users becomes the active/selected sheet with selection D9
$usersWs.Activate() | Out-Null
$usersWs.Range("D9").Select() | Out-Null
